$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New period order (most recent first) with each period's original "Valor Mora" value.
# The previous statements (periods 1804-2412) are superseded; new/reversed period list populates B16:J98.
$periodos = @("2502","2501","2412","2411","2410","2409","2408","2407","2406","2405","2404","2403","2402","2401","2312","2311","2310","2309","2308","2307","2306","2305","2304","2303","2302","2301","2212","2211","2210","2209","2208","2207","2206","2205","2204","2203","2202","2201","2112","2111","2110","2109","2108","2107","2106","2105","2104","2103","2102","2101","2012","2011","2010","2009","2008","2007","2006","2005","2004","2003","2002","2001","1912","1911","1910","1909","1908","1907","1906","1905","1904","1903","1902","1901","1812","1811","1810","1809","1808","1807","1806","1805","1804")
$valores = @(19791,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,31249,7812,7812,7812,7812,7812,7812,7812,7812,7812,7812,7812,7812,7812,5208)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
